# Update crypto price/volume symbol list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values must remain plain text (matching the
# workbook's existing inline-string cells), so we force the Text number
# format before assigning, which preserves exact formatting (trailing
# zeros, '%' suffix, etc.) instead of letting Excel auto-coerce to a number.
$changes = @(
    @{ Cell = "D2"; Value = "314.75" },
    @{ Cell = "E2"; Value = "3.39%" },
    @{ Cell = "D3"; Value = "35.64" },
    @{ Cell = "E3"; Value = "0.09%" },
    @{ Cell = "D4"; Value = "5.098" },
    @{ Cell = "E4"; Value = "0.90%" },
    @{ Cell = "D5"; Value = "0.08151" },
    @{ Cell = "E5"; Value = "3.29%" },
    @{ Cell = "D6"; Value = "2.105" },
    @{ Cell = "E6"; Value = "-0.62%" },
    @{ Cell = "D7"; Value = "4.143" },
    @{ Cell = "E7"; Value = "0.53%" },
    @{ Cell = "D8"; Value = "7.938" },
    @{ Cell = "E8"; Value = "0.41%" },
    @{ Cell = "D9"; Value = "0.9335" },
    @{ Cell = "E9"; Value = "1.10%" },
    @{ Cell = "D10"; Value = "0.1036" },
    @{ Cell = "E10"; Value = "7.23%" },
    @{ Cell = "D11"; Value = "0.1921" },
    @{ Cell = "E11"; Value = "4.30%" },
    @{ Cell = "D12"; Value = "0.09025" },
    @{ Cell = "E12"; Value = "4.25%" },
    @{ Cell = "D13"; Value = "0.03614" },
    @{ Cell = "E13"; Value = "1.58%" },
    @{ Cell = "D14"; Value = "0.09885" },
    @{ Cell = "E14"; Value = "-0.14%" },
    @{ Cell = "D15"; Value = "0.001430" },
    @{ Cell = "E15"; Value = "-0.33%" },
    @{ Cell = "D16"; Value = "0.005838" },
    @{ Cell = "E16"; Value = "2.99%" },
    @{ Cell = "E17"; Value = "-0.22%" },
    @{ Cell = "D18"; Value = "2.977" },
    @{ Cell = "E18"; Value = "12.76%" },
    @{ Cell = "D19"; Value = "0.3459" },
    @{ Cell = "E19"; Value = "2.64%" },
    @{ Cell = "D20"; Value = "0.1313" },
    @{ Cell = "E20"; Value = "-1.89%" },
    @{ Cell = "D21"; Value = "5.118" },
    @{ Cell = "E21"; Value = "-0.79%" },
    @{ Cell = "E22"; Value = "0.32%" },
    @{ Cell = "D23"; Value = "0.04551" },
    @{ Cell = "E23"; Value = "1.05%" },
    @{ Cell = "D24"; Value = "0.001244" },
    @{ Cell = "E24"; Value = "0.92%" },
    @{ Cell = "E25"; Value = "-1.34%" },
    @{ Cell = "D26"; Value = "0.0001252" },
    @{ Cell = "E26"; Value = "-3.84%" },
    @{ Cell = "E27"; Value = "-5.11%" },
    @{ Cell = "D39"; Value = "0.01961" },
    @{ Cell = "E39"; Value = "6.64%" },
    @{ Cell = "D40"; Value = "0.04898" },
    @{ Cell = "E40"; Value = "3.40%" },
    @{ Cell = "D41"; Value = "0.007590" },
    @{ Cell = "E41"; Value = "-3.63%" },
    @{ Cell = "E42"; Value = "-0.46%" },
    @{ Cell = "D43"; Value = "0.007890" },
    @{ Cell = "E43"; Value = "1.63%" },
    @{ Cell = "D44"; Value = "0.002098" },
    @{ Cell = "E44"; Value = "-4.29%" },
    @{ Cell = "D45"; Value = "0.01176" },
    @{ Cell = "E45"; Value = "5.45%" },
    @{ Cell = "D46"; Value = "0.00006737" },
    @{ Cell = "E46"; Value = "6.92%" },
    @{ Cell = "E47"; Value = "0.15%" },
    @{ Cell = "D48"; Value = "186.75" },
    @{ Cell = "E48"; Value = "268.96%" },
    @{ Cell = "E49"; Value = "-10.40%" },
    @{ Cell = "D50"; Value = "0.00002106" },
    @{ Cell = "E50"; Value = "0.15%" },
    @{ Cell = "D51"; Value = "0.0002006" },
    @{ Cell = "E51"; Value = "0.15%" }
)

foreach ($chg in $changes) {
    $cell = $ws.Range($chg.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
}
